$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force a text number-format on the cells whose new values look like
# numbers, so Excel stores them as text (matching the scraped source
# data, e.g. "25.738.10") instead of silently converting to a Double.
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D51').NumberFormat = '@'

$ws.Range('D2').Value = '25.738.10'
$ws.Range('E2').Value = '  -2.88%  '
$ws.Range('D3').Value = '1.744.34'
$ws.Range('E3').Value = '  -5.01%  '
$ws.Range('E4').Value = '  -0.03%  '
$ws.Range('D5').Value = '237.68'
$ws.Range('E5').Value = '  -7.57%  '
$ws.Range('E6').Value = '  -0.04%  '
$ws.Range('D7').Value = '0.4901'
$ws.Range('E7').Value = '  -6.53%  '
$ws.Range('D8').Value = '41.76'
$ws.Range('E8').Value = '  -6.83%  '
$ws.Range('D9').Value = '0.2423'
$ws.Range('E9').Value = '  -22.98%  '
$ws.Range('D10').Value = '0.06027'
$ws.Range('E10').Value = '  -11.03%  '
$ws.Range('D11').Value = '1.761.47'
$ws.Range('E11').Value = '  -3.93%  '
$ws.Range('D12').Value = '0.06662'
$ws.Range('E12').Value = '  -14.13%  '
$ws.Range('D13').Value = '14.47'
$ws.Range('E13').Value = '  -22.15%  '
$ws.Range('D14').Value = '0.5978'
$ws.Range('E14').Value = '  -22.77%  '
$ws.Range('D15').Value = '77.25'
$ws.Range('E15').Value = '  -11.72%  '
$ws.Range('D16').Value = '4.325'
$ws.Range('E16').Value = '  -13.34%  '
$ws.Range('D17').Value = '1.001'
$ws.Range('E17').Value = '  -0.04%  '
$ws.Range('E18').Value = '  -0.06%  '
$ws.Range('D19').Value = '25.770.32'
$ws.Range('E19').Value = '  -2.85%  '
$ws.Range('D20').Value = '11.24'
$ws.Range('E20').Value = '  -18.45%  '
$ws.Range('D21').Value = '0.000006275'
$ws.Range('E21').Value = '  -20.57%  '
$ws.Range('D22').Value = '1.979.65'
$ws.Range('E22').Value = '  -4.28%  '
$ws.Range('D23').Value = '3.869'
$ws.Range('E23').Value = '  -15.61%  '
$ws.Range('D24').Value = '5.087'
$ws.Range('E24').Value = '  -14.47%  '
$ws.Range('D25').Value = '7.961'
$ws.Range('E25').Value = '  -14.16%  '
$ws.Range('D26').Value = '134.41'
$ws.Range('D27').Value = '1.499'
$ws.Range('E27').Value = '  -10.37%  '
$ws.Range('D28').Value = '1.864'
$ws.Range('E28').Value = '  -15.28%  '
$ws.Range('D29').Value = '14.32'
$ws.Range('E29').Value = '  -15.18%  '
$ws.Range('D30').Value = '99.18'
$ws.Range('E30').Value = '  -10.74%  '
$ws.Range('D31').Value = '0.08180'
$ws.Range('E31').Value = '  -6.18%  '
$ws.Range('D32').Value = '3.618'
$ws.Range('E32').Value = '  -12.68%  '
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').Value = '3.169'
$ws.Range('E33').Value = '  -21.69%  '
$ws.Range('B34').Value = 'Hedera'
$ws.Range('C34').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D34').Value = '0.04291'
$ws.Range('E34').Value = '  -11.44%  '
$ws.Range('B35').Value = 'ARBITRUM'
$ws.Range('C35').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D35').Value = '1.032'
$ws.Range('E35').Value = '  -8.74%  '
$ws.Range('B36').Value = 'HuobiToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D36').Value = '2.611'
$ws.Range('E36').Value = '  -8.79%  '
$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D37').Value = '0.6073'
$ws.Range('E37').Value = '  -15.07%  '
$ws.Range('B38').Value = 'MXToken'
$ws.Range('C38').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range('D38').Value = '2.770'
$ws.Range('E38').Value = '  -10.22%  '
$ws.Range('B39').Value = 'RenderToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D39').Value = '2.080'
$ws.Range('E39').Value = '  -6.00%  '
$ws.Range('B40').Value = 'PaxDollar'
$ws.Range('C40').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D40').Value = '1.001'
$ws.Range('E40').Value = '  -0.08%  '
$ws.Range('B41').Value = 'VeChain'
$ws.Range('C41').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D41').Value = '0.01481'
$ws.Range('E41').Value = '  -14.18%  '
$ws.Range('D42').Value = '101.02'
$ws.Range('E42').Value = '  -7.85%  '
$ws.Range('B43').Value = 'TrustWalletToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D43').Value = '0.7886'
$ws.Range('E43').Value = '  -11.31%  '
$ws.Range('B44').Value = 'TheSandbox'
$ws.Range('C44').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D44').Value = '0.3816'
$ws.Range('E44').Value = '  -20.23%  '
$ws.Range('B45').Value = 'FraxShare'
$ws.Range('C45').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D45').Value = '5.126'
$ws.Range('E45').Value = '  -13.31%  '
$ws.Range('B46').Value = 'Aptos'
$ws.Range('C46').Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range('D46').Value = '6.076'
$ws.Range('E46').Value = '  -20.05%  '
$ws.Range('B47').Value = 'Cronos'
$ws.Range('C47').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D47').Value = '0.05073'
$ws.Range('E47').Value = '  -12.57%  '
$ws.Range('B48').Value = 'Aave'
$ws.Range('C48').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D48').Value = '51.85'
$ws.Range('E48').Value = '  -12.71%  '
$ws.Range('D49').Value = '29.51'
$ws.Range('E49').Value = '  -14.78%  '
$ws.Range('B50').Value = 'Algorand'
$ws.Range('C50').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D50').Value = '0.1032'
$ws.Range('E50').Value = '  -15.56%  '
$ws.Range('B51').Value = 'USDD'
$ws.Range('C51').Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range('D51').Value = '0.9991'
$ws.Range('E51').Value = '  -0.42%  '

# Clear the temporary text formatting so these cells fall back to the
# default style (no explicit "s" attribute), matching the original
# workbook formatting for these data cells.
$ws.Range('D5').ClearFormats()
$ws.Range('D7').ClearFormats()
$ws.Range('D8').ClearFormats()
$ws.Range('D9').ClearFormats()
$ws.Range('D10').ClearFormats()
$ws.Range('D12').ClearFormats()
$ws.Range('D13').ClearFormats()
$ws.Range('D14').ClearFormats()
$ws.Range('D15').ClearFormats()
$ws.Range('D16').ClearFormats()
$ws.Range('D17').ClearFormats()
$ws.Range('D20').ClearFormats()
$ws.Range('D21').ClearFormats()
$ws.Range('D23').ClearFormats()
$ws.Range('D24').ClearFormats()
$ws.Range('D25').ClearFormats()
$ws.Range('D26').ClearFormats()
$ws.Range('D27').ClearFormats()
$ws.Range('D28').ClearFormats()
$ws.Range('D29').ClearFormats()
$ws.Range('D30').ClearFormats()
$ws.Range('D31').ClearFormats()
$ws.Range('D32').ClearFormats()
$ws.Range('D33').ClearFormats()
$ws.Range('D34').ClearFormats()
$ws.Range('D35').ClearFormats()
$ws.Range('D36').ClearFormats()
$ws.Range('D37').ClearFormats()
$ws.Range('D38').ClearFormats()
$ws.Range('D39').ClearFormats()
$ws.Range('D40').ClearFormats()
$ws.Range('D41').ClearFormats()
$ws.Range('D42').ClearFormats()
$ws.Range('D43').ClearFormats()
$ws.Range('D44').ClearFormats()
$ws.Range('D45').ClearFormats()
$ws.Range('D46').ClearFormats()
$ws.Range('D47').ClearFormats()
$ws.Range('D48').ClearFormats()
$ws.Range('D49').ClearFormats()
$ws.Range('D50').ClearFormats()
$ws.Range('D51').ClearFormats()
